# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C for rows 2-8 by one day
# (date serial 46059 -> 46060, i.e. 2026-02-06 -> 2026-02-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46059) {
        $cell.Value = 46060
    }
}
